# Adicion de carpeta web: se agrega una columna "url" calculada en cada
# hoja (obras / referentes) que concatena la ruta publica del repositorio
# con el nombre de archivo existente.

$wb = $excel.ActiveWorkbook

$wsObras = $wb.Worksheets.Item("obras")
$wsReferentes = $wb.Worksheets.Item("referentes")

# --- Hoja "referentes" (se procesa primero para que la hoja que quede
#     activa al final, tal como en el archivo original, sea "obras") -----
$wsReferentes.Activate()

$wsReferentes.Range("E4").Value = "laminas-paisajes.jpg"

$wsReferentes.Range("F1").Value = "url"
$wsReferentes.Range("F1").Style = $wsReferentes.Range("E1").Style
$wsReferentes.Range("F2").Formula = '=CONCATENATE("https://badac.uniandes.edu.co/files/bga-cat/",E2)'
$wsReferentes.Range("F3").Formula = '=CONCATENATE("https://badac.uniandes.edu.co/files/bga-cat/",E3)'
$wsReferentes.Range("F4").Formula = '=CONCATENATE("https://badac.uniandes.edu.co/files/bga-cat/",E4)'
$wsReferentes.Range("F5").Formula = '=CONCATENATE("https://badac.uniandes.edu.co/files/bga-cat/",E5)'

$wsReferentes.Columns.Item(6).ColumnWidth = 69.7109375

$wsReferentes.Range("H12").Font.Underline = $true

$wsReferentes.Range("F3").Select()
$wsReferentes.Application.ActiveWindow.ScrollColumn = 3

# --- Hoja "obras" --------------------------------------------------------
# La columna "Archivo" vive en H, pero la formula original (tal y como fue
# escrita por la autora) referencia la columna E (Ancho cm) -- se respeta
# tal cual para reproducir el archivo real.
$wsObras.Activate()

$wsObras.Range("J1").Value = "url"
$wsObras.Range("J1").Style = $wsObras.Range("I1").Style
$wsObras.Range("J2").Formula = '=CONCATENATE("https://badac.uniandes.edu.co/files/bga-cat/",E2)'
$wsObras.Range("J3").Formula = '=CONCATENATE("https://badac.uniandes.edu.co/files/bga-cat/",E3)'
$wsObras.Range("J4").Formula = '=CONCATENATE("https://badac.uniandes.edu.co/files/bga-cat/",E4)'

$wsObras.Columns.Item(10).ColumnWidth = 47

# celda suelta con formato (subrayado) que quedo marcada en la hoja tras
# el taller
$wsObras.Range("L11").Font.Underline = $true

$wsObras.Range("L11").Select()
$wsObras.Application.ActiveWindow.ScrollColumn = 8
